$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values look numeric,
# so Excel keeps them as literal text instead of auto-converting to a number.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the updated cell values
$ws.Range("D2").Value = '29.719.72'
$ws.Range("E2").Value = '  -1.11%  '
$ws.Range("D3").Value = '1.886.24'
$ws.Range("E3").Value = '  -1.33%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '0.7915'
$ws.Range("E5").Value = '  -4.38%  '
$ws.Range("D6").Value = '241.18'
$ws.Range("E6").Value = '  -0.48%  '
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("D8").Value = '0.3163'
$ws.Range("E8").Value = '  -2.07%  '
$ws.Range("D9").Value = '25.46'
$ws.Range("E9").Value = '  -4.87%  '
$ws.Range("D10").Value = '0.06989'
$ws.Range("D11").Value = '0.08040'
$ws.Range("E11").Value = '  +0.03%  '
$ws.Range("D12").Value = '0.7571'
$ws.Range("E12").Value = '  +0.78%  '
$ws.Range("D13").Value = '1.903.45'
$ws.Range("E13").Value = '  -0.43%  '
$ws.Range("D14").Value = '5.273'
$ws.Range("E14").Value = '  +0.84%  '
$ws.Range("D15").Value = '92.04'
$ws.Range("D16").Value = '29.755.66'
$ws.Range("E16").Value = '  -1.01%  '
$ws.Range("D17").Value = '13.81'
$ws.Range("E17").Value = '  -2.70%  '
$ws.Range("D18").Value = '5.916'
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("D19").Value = '243.03'
$ws.Range("E19").Value = '  -0.85%  '
$ws.Range("D20").Value = '0.000007667'
$ws.Range("E20").Value = '  -1.40%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.16%  '
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").Value = '2.150.14'
$ws.Range("E22").Value = '  -0.54%  '
$ws.Range("D23").Value = '8.158'
$ws.Range("E23").Value = '  +16.67%  '
$ws.Range("D24").Value = '1.001'
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("D25").Value = '0.1660'
$ws.Range("D26").Value = '9.284'
$ws.Range("E26").Value = '  +0.19%  '
$ws.Range("D27").Value = '164.13'
$ws.Range("E27").Value = '  -3.42%  '
$ws.Range("E28").Value = '  -1.96%  '
$ws.Range("D29").Value = '2.042'
$ws.Range("E29").Value = '  -2.14%  '
$ws.Range("D30").Value = '1.397'
$ws.Range("E30").Value = '  +2.02%  '
$ws.Range("D31").Value = '1.531'
$ws.Range("E31").Value = '  +0.76%  '
$ws.Range("D32").Value = '4.373'
$ws.Range("E32").Value = '  +1.47%  '
$ws.Range("D33").Value = '0.05672'
$ws.Range("E33").Value = '  +1.27%  '
$ws.Range("D34").Value = '4.043'
$ws.Range("E34").Value = '  -1.40%  '
$ws.Range("E35").Value = '  -2.20%  '
$ws.Range("D36").Value = '0.7319'
$ws.Range("E36").Value = '  -0.64%  '
$ws.Range("D37").Value = '0.9968'
$ws.Range("E37").Value = '  -0.18%  '
$ws.Range("D38").Value = '2.607'
$ws.Range("E38").Value = '  -4.12%  '
$ws.Range("D39").Value = '0.01903'
$ws.Range("E39").Value = '  -0.89%  '
$ws.Range("D40").Value = '2.768'
$ws.Range("E40").Value = '  -1.00%  '
$ws.Range("D41").Value = '0.4391'
$ws.Range("E41").Value = '  -1.16%  '
$ws.Range("D42").Value = '72.31'
$ws.Range("E42").Value = '  -0.41%  '
$ws.Range("D43").Value = '5.799'
$ws.Range("E43").Value = '  -3.33%  '
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.16%  '
$ws.Range("D45").Value = '0.8355'
$ws.Range("E45").Value = '  -0.82%  '
$ws.Range("D46").Value = '102.57'
$ws.Range("E46").Value = '  +1.45%  '
$ws.Range("D47").Value = '1.019.07'
$ws.Range("E47").Value = '  +3.59%  '
$ws.Range("D48").Value = '1.860'
$ws.Range("E48").Value = '  -2.24%  '
$ws.Range("D49").Value = '9.878'
$ws.Range("E49").Value = '  +0.91%  '
$ws.Range("D50").Value = '7.412'
$ws.Range("E50").Value = '  -2.74%  '
$ws.Range("D51").Value = '2.040.89'
$ws.Range("E51").Value = '  -1.08%  '
